$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "HK_G_acc_LG"

$values = @(
    95.439189189189193,
    96.11486486486487,
    95.777027027027032,
    95.270270270270274,
    95.608108108108098,
    95.945945945945937,
    94.763513513513516,
    95.608108108108098,
    95.608108108108098,
    95.101351351351354,
    88.344594594594597,
    94.594594594594597,
    94.087837837837839,
    93.074324324324323,
    93.074324324324323,
    91.722972972972968,
    94.425675675675677,
    94.425675675675677,
    93.243243243243242,
    95.101351351351354,
    94.087837837837839,
    91.21621621621621,
    92.567567567567565,
    92.567567567567565,
    96.28378378378379,
    94.932432432432435,
    96.11486486486487,
    93.581081081081081,
    94.087837837837839,
    95.439189189189193,
    95.945945945945937,
    95.608108108108098,
    95.608108108108098,
    94.087837837837839,
    94.594594594594597,
    71.452702702702695,
    94.256756756756758,
    94.932432432432435,
    93.243243243243242,
    95.777027027027032,
    95.270270270270274,
    95.777027027027032,
    95.270270270270274,
    94.932432432432435,
    95.270270270270274,
    95.270270270270274,
    89.189189189189193,
    95.777027027027032
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
